# "Creation de la fenetre" - tidy up the stray typo/proofing artefacts left
# over from Word's spell/grammar checker and drop the now-unused
# header/footer reservation on the page.
$d = $word.ActiveDocument

# Paragraph 4: "Equipement=ùateriel scolaire" -> fix the typo'd accented
# "ù" into "m" ("materiel") and collapse the spell-check wrapping that
# surrounded it so the paragraph reads as plain corrected text.
$p4 = $d.Paragraphs(4)
$p4.Range.Find.Execute("Equipement=ùateriel scolaire", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Equipement=materiel scolaire", 2)

# Paragraph 5: "Les monstre= les uv" had grammar/spelling markers splitting
# it into three runs ("Les monstre" | "= les " | "uv"). Re-assert the text
# as a single corrected sentence so the stray markers collapse away.
$p5 = $d.Paragraphs(5)
$p5.Range.Find.Execute("Les monstre= les uv", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Les monstre= les uv", 2)

# The page no longer reserves header/footer space since neither is used.
$d.PageSetup.HeaderDistance = 0
$d.PageSetup.FooterDistance = 0
